$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new row 53 gets the same style as the date column above (s="2")
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 7.226520411029047
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 8.987952903094421
$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 6.109216616889168
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 8.045645122021906
$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 2008
$ws.Cells.Item(4, 3).Value = 4.268860212333636
$ws.Cells.Item(4, 4).Value = 2009
$ws.Cells.Item(4, 5).Value = 6.77211531652997
$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = -5.232639093663815
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = -0.9943400519801915
$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2009
$ws.Cells.Item(6, 3).Value = -7.266312015249776
$ws.Cells.Item(6, 4).Value = 2010
$ws.Cells.Item(6, 5).Value = -1.420242831007679
$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 2.12454018480297
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = 1.203634802640963
$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 2010
$ws.Cells.Item(8, 3).Value = 6.958243460951929
$ws.Cells.Item(8, 4).Value = 2011
$ws.Cells.Item(8, 5).Value = 8.31992099499319
$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 8.081020954067775
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = 2.257871268432821
$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 2011
$ws.Cells.Item(10, 3).Value = 9.469137444079934
$ws.Cells.Item(10, 4).Value = 2012
$ws.Cells.Item(10, 5).Value = 8.571528775834981
$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 4.489210662380949
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 7.472658273721078
$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 2012
$ws.Cells.Item(12, 3).Value = 3.358206407534947
$ws.Cells.Item(12, 4).Value = 2013
$ws.Cells.Item(12, 5).Value = 5.745831525574441
$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = -0.8752093743685352
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 1.998870338019265
$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 2013
$ws.Cells.Item(14, 3).Value = 0.3081076735359067
$ws.Cells.Item(14, 4).Value = 2014
$ws.Cells.Item(14, 5).Value = 2.847379875994704
$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 4.891728508251214
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 2.926340920335191
$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 2014
$ws.Cells.Item(16, 3).Value = 3.901355411819707
$ws.Cells.Item(16, 4).Value = 2015
$ws.Cells.Item(16, 5).Value = 3.690459963535009
$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = 4.818339085077583
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 4.241902819910548
$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 5.331683351557981
$ws.Cells.Item(18, 4).Value = 2016
$ws.Cells.Item(18, 5).Value = 4.636575318346536
$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 4.067959312311897
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 4.36243732366437
$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 2016
$ws.Cells.Item(20, 3).Value = 3.254758369308375
$ws.Cells.Item(20, 4).Value = 2017
$ws.Cells.Item(20, 5).Value = 2.76788332063731
$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 4.613634856640747
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 4.246555741688218
$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = 5.246209615995667
$ws.Cells.Item(22, 4).Value = 2018
$ws.Cells.Item(22, 5).Value = 4.659862065074982
$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 6.155351106582874
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 5.081201994458495
$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 3.898744563937395
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 3.529300656691237
$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 5.021833280292598
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 4.83330568333058
$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 2018
$ws.Cells.Item(26, 3).Value = 4.86255966374296
$ws.Cells.Item(26, 4).Value = 2019
$ws.Cells.Item(26, 5).Value = 4.636196713604357
$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 4.237941638782527
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 4.932073907517265
$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 3.841510956591465
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 4.83848589746565
$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 3.042742667481324
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 3.228901920070748
$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = 2.764740011159428
$ws.Cells.Item(30, 4).Value = 2020
$ws.Cells.Item(30, 5).Value = 2.471557257221946
$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 2.059284312217757
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 2.845541644111571
$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = -0.8225206269755425
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = 1.159653508089242
$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -7.578477024949737
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = -5.048497756254311
$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = -7.260793671746435
$ws.Cells.Item(34, 4).Value = 2021
$ws.Cells.Item(34, 5).Value = -1.387795042833839
$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = -0.4318290737559183
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = -1.967992356539539
$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = 4.409066926520455
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = 3.034999751677669
$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = 4.379227219808146
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = 2.570052815033752
$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 2021
$ws.Cells.Item(38, 3).Value = 4.097586525396268
$ws.Cells.Item(38, 4).Value = 2022
$ws.Cells.Item(38, 5).Value = 2.450242954096926
$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 6.557154773664364
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = 1.465110550383386
$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 8.064077385547574
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = 3.690055931494096
$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 7.397318165265498
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = 3.192017772210276
$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 2022
$ws.Cells.Item(42, 3).Value = 7.824284864703746
$ws.Cells.Item(42, 4).Value = 2023
$ws.Cells.Item(42, 5).Value = 4.834990656989402
$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = 1.287151040638124
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = 5.402455143891816
$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = 0.3551698673347259
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = 4.038118345571751
$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = -0.3046246622258053
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = 2.185848087053199
$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 2023
$ws.Cells.Item(46, 3).Value = -1.24502235313334
$ws.Cells.Item(46, 4).Value = 2024
$ws.Cells.Item(46, 5).Value = -2.488220481262082
$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = -3.197915208378399
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = 0.5664860188349996
$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = -2.267078452724969
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 0.5378929214800987
$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = -2.567041707495976
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = -0.3140572462435154
$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 2024
$ws.Cells.Item(50, 3).Value = -1.735114423676209
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = 1.056286187957367
$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = 1.602010908728335
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = -1.246071640539481
$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = 1.326993065386817
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = -0.2087957186147071
$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 2025
$ws.Cells.Item(53, 3).Value = 2.64031107104763
$ws.Cells.Item(53, 4).Value = 2026
$ws.Cells.Item(53, 5).Value = 1.890175284700679
